$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.980.59'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '2.647.83'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.621'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.88%  '
$ws.Range('D9').Value = '2.646.32'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.118'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.81'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.63'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').Value = '3.123.58'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000185'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '63.918.35'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '2.669.98'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.79%  '
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '346.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +5.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000112'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '585.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.19%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.07'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.64'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.50%  '
$ws.Range('E37').Value = '  -2.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '151.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.18%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.55'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.52%  '
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '162.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.33%  '
$ws.Range('E47').Value = '  -2.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0591'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.635'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0250'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.100'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.03%  '
